# Continuing Issue686 [Update Features.html to match look and feel of the
# home page] - Update Issue 686 Some more tweaks.
#
# The reviewer comment left on slide 1 ("Button not visible. Not aligned
# with anything.") is no longer relevant, so remove it. This drops
# ppt/comments/comment1.xml (and the now-dangling content-type override /
# slide relationship that point at it) while leaving the comment author
# list (ppt/commentAuthors.xml) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Comments.Count; $i -ge 1; $i--) {
    $s.Comments.Item($i).Delete()
}
